$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 975.3333
$ws.Range("I96").Value = 227.88889
$ws.Range("J96").Value = 2096.5
$ws.Range("K96").Value = 683.6666700000001
$ws.Range("L96").Value = 6289.5
$ws.Range("M96").Value = 689.3333299999999
$ws.Range("N96").Value = -9035.5

$ws.Range("H125").Value = 1000.25
$ws.Range("I125").Value = 900
$ws.Range("J125").Value = 1100.5
$ws.Range("K125").Value = 8100
$ws.Range("L125").Value = 9904.5
$ws.Range("M125").Value = -5640
$ws.Range("N125").Value = -14824.5

$ws.Range("H129").Value = 932.0769
$ws.Range("J129").Value = 1068.2858
$ws.Range("L129").Value = 3204.8574
$ws.Range("N129").Value = -13204.8574

$ws.Range("H138").Value = 3178.2727
$ws.Range("I138").Value = 2485.2856
$ws.Range("J138").Value = 3309.3784
$ws.Range("K138").Value = 7455.8568
$ws.Range("L138").Value = 9928.135200000001
$ws.Range("M138").Value = -2315.8568
$ws.Range("N138").Value = -20208.1352

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16963.955
$ws.Range("I32").Value = 4057.6
$ws.Range("J32").Value = 42776.668
$ws.Range("K32").Value = 4057.6
$ws.Range("L32").Value = 42776.668
$ws.Range("M32").Value = -3770.6
$ws.Range("N32").Value = -43350.668

$ws.Range("H97").Value = 2880.0715
$ws.Range("I97").Value = 3245.0833
$ws.Range("J97").Value = 690
$ws.Range("K97").Value = 3245.0833
$ws.Range("L97").Value = 690
$ws.Range("M97").Value = -2749.0833
$ws.Range("N97").Value = -1682

$ws.Range("H102").Value = 1701.25
$ws.Range("I102").Value = 1801.4286
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1801.4286
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -179.4286
$ws.Range("N102").Value = -4244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 50002244
$ws.Range("I86").Value = 50002244
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 50002244
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -50001121
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 50002244
$ws.Range("I89").Value = 50002244
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 250011220
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -250005604
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 7306.645
$ws.Range("I94").Value = 858.2857
$ws.Range("J94").Value = 20848.2
$ws.Range("K94").Value = 858.2857
$ws.Range("L94").Value = 20848.2
$ws.Range("M94").Value = -407.2857
$ws.Range("N94").Value = -21750.2

$ws.Range("H99").Value = 1705.2354
$ws.Range("I99").Value = 1784.4286
$ws.Range("J99").Value = 1649.8
$ws.Range("K99").Value = 1784.4286
$ws.Range("L99").Value = 1649.8
$ws.Range("M99").Value = -286.4286
$ws.Range("N99").Value = -4645.8

$ws.Range("H105").Value = 6300
$ws.Range("I105").Value = 2271.4285
$ws.Range("J105").Value = 11940
$ws.Range("K105").Value = 2271.4285
$ws.Range("L105").Value = 11940
$ws.Range("M105").Value = -524.4285
$ws.Range("N105").Value = -15434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3392286
$ws.Range("I31").Value = 10528306
$ws.Range("J31").Value = 2676.7
$ws.Range("K31").Value = 10528306
$ws.Range("L31").Value = 2676.7
$ws.Range("M31").Value = -10528011
$ws.Range("N31").Value = -3266.7

$ws.Range("H34").Value = 3392286
$ws.Range("I34").Value = 10528306
$ws.Range("J34").Value = 2676.7
$ws.Range("K34").Value = 10528306
$ws.Range("L34").Value = 2676.7
$ws.Range("M34").Value = -10528104
$ws.Range("N34").Value = -3080.7

$ws.Range("H62").Value = 3106.1765
$ws.Range("I62").Value = 2734.1667
$ws.Range("J62").Value = 3309.0908
$ws.Range("K62").Value = 2734.1667
$ws.Range("L62").Value = 3309.0908
$ws.Range("M62").Value = -2110.1667
$ws.Range("N62").Value = -4557.0908

$ws.Range("H65").Value = 3106.1765
$ws.Range("I65").Value = 2734.1667
$ws.Range("J65").Value = 3309.0908
$ws.Range("K65").Value = 13670.8335
$ws.Range("L65").Value = 16545.454
$ws.Range("M65").Value = -10550.8335
$ws.Range("N65").Value = -22785.454

$ws.Range("H105").Value = 850
$ws.Range("I105").Value = 775
$ws.Range("K105").Value = 775
$ws.Range("M105").Value = 972

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 6382.606
$ws.Range("I87").Value = 4350
$ws.Range("J87").Value = 6513.7417
$ws.Range("K87").Value = 13050
$ws.Range("L87").Value = 19541.2251
$ws.Range("M87").Value = -11802
$ws.Range("N87").Value = -22037.2251

$ws.Range("H90").Value = 6382.606
$ws.Range("I90").Value = 4350
$ws.Range("J90").Value = 6513.7417
$ws.Range("K90").Value = 39150
$ws.Range("L90").Value = 58623.6753
$ws.Range("M90").Value = -32910
$ws.Range("N90").Value = -71103.6753

$ws.Range("H131").Value = 1961817.2
$ws.Range("I131").Value = 9524191
$ws.Range("J131").Value = 1201.963
$ws.Range("K131").Value = 28572573
$ws.Range("L131").Value = 3605.889
$ws.Range("M131").Value = -28567533
$ws.Range("N131").Value = -13685.889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 3009
$ws.Range("J25").Value = 3009
$ws.Range("L25").Value = 3009
$ws.Range("N25").Value = -4067

$ws.Range("H33").Value = 9754.75
$ws.Range("J33").Value = 9754.75
$ws.Range("L33").Value = 9754.75
$ws.Range("N33").Value = -10258.75

$ws.Range("H80").Value = 3985
$ws.Range("I80").Value = 7276.25
$ws.Range("J80").Value = 2522.2222
$ws.Range("K80").Value = 7276.25
$ws.Range("L80").Value = 2522.2222
$ws.Range("M80").Value = -6278.25
$ws.Range("N80").Value = -4518.2222

$ws.Range("H83").Value = 3985
$ws.Range("I83").Value = 7276.25
$ws.Range("J83").Value = 2522.2222
$ws.Range("K83").Value = 36381.25
$ws.Range("L83").Value = 12611.111
$ws.Range("M83").Value = -31389.25
$ws.Range("N83").Value = -22595.111

$ws.Range("H97").Value = 1167.9286
$ws.Range("I97").Value = 1208.2727
$ws.Range("J97").Value = 1020
$ws.Range("K97").Value = 1208.2727
$ws.Range("L97").Value = 1020
$ws.Range("M97").Value = -712.2727
$ws.Range("N97").Value = -2012

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1930.2858
$ws.Range("I7").Value = 1197.3334
$ws.Range("J7").Value = 2480
$ws.Range("K7").Value = 1197.3334
$ws.Range("L7").Value = 2480
$ws.Range("M7").Value = -1085.3334
$ws.Range("N7").Value = -2704

$ws.Range("H68").Value = 19904294
$ws.Range("I68").Value = 112778250
$ws.Range("J68").Value = 2731.7144
$ws.Range("K68").Value = 112778250
$ws.Range("L68").Value = 2731.7144
$ws.Range("M68").Value = -112777501
$ws.Range("N68").Value = -4229.7144

$ws.Range("H71").Value = 19904294
$ws.Range("I71").Value = 112778250
$ws.Range("J71").Value = 2731.7144
$ws.Range("K71").Value = 563891250
$ws.Range("L71").Value = 13658.572
$ws.Range("M71").Value = -563887506
$ws.Range("N71").Value = -21146.572

$ws.Range("H82").Value = 1841
$ws.Range("I82").Value = 1536.6666
$ws.Range("J82").Value = 1971.4286
$ws.Range("K82").Value = 1536.6666
$ws.Range("L82").Value = 1971.4286
$ws.Range("M82").Value = -1175.6666
$ws.Range("N82").Value = -2693.4286

$ws.Range("H85").Value = 1841
$ws.Range("I85").Value = 1536.6666
$ws.Range("J85").Value = 1971.4286
$ws.Range("K85").Value = 1536.6666
$ws.Range("L85").Value = 1971.4286
$ws.Range("M85").Value = -288.6666
$ws.Range("N85").Value = -4467.4286

$ws.Range("H93").Value = 2028.0834
$ws.Range("I93").Value = 1722.4286
$ws.Range("J93").Value = 2456
$ws.Range("K93").Value = 1722.4286
$ws.Range("L93").Value = 2456
$ws.Range("M93").Value = -474.4286
$ws.Range("N93").Value = -4952

$ws.Range("H100").Value = 1587.875
$ws.Range("I100").Value = 1483.8334
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 1483.8334
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -942.8334
$ws.Range("N100").Value = -2982

$ws.Range("H122").Value = 4111.6665
$ws.Range("I122").Value = 3250
$ws.Range("K122").Value = 9750
$ws.Range("M122").Value = -7300

$ws.Range("H126").Value = 1930.2858
$ws.Range("I126").Value = 1197.3334
$ws.Range("J126").Value = 2480
$ws.Range("K126").Value = 3592.0002
$ws.Range("L126").Value = 7440
$ws.Range("M126").Value = -1122.0002
$ws.Range("N126").Value = -12380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 66669188
$ws.Range("I81").Value = 125002550
$ws.Range("J81").Value = 2485.7144
$ws.Range("K81").Value = 250005100
$ws.Range("L81").Value = 4971.4288
$ws.Range("M81").Value = -250004039
$ws.Range("N81").Value = -7093.4288

$ws.Range("H84").Value = 66669188
$ws.Range("I84").Value = 125002550
$ws.Range("J84").Value = 2485.7144
$ws.Range("K84").Value = 1250025500
$ws.Range("L84").Value = 24857.144
$ws.Range("M84").Value = -1250020196
$ws.Range("N84").Value = -35465.144

$ws.Range("H96").Value = 1200
$ws.Range("I96").Value = 1200
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 173
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 1083.25
$ws.Range("I122").Value = 894.5789
$ws.Range("J122").Value = 1359
$ws.Range("K122").Value = 2683.7367
$ws.Range("L122").Value = 4077
$ws.Range("M122").Value = -233.7366999999999
$ws.Range("N122").Value = -8977

Write-Host "Applied scheduled market-data update across all 8 sheets"
